$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 98.23421133913619
$ws.Range("E2").Value = 105.0149002075195
$ws.Range("F2").Value = 107.3691364540896
$ws.Range("G2").Value = 97.10613185382886
$ws.Range("H2").Value = 934735206
$ws.Range("I2").Value = "IBM"

$ws.Range("D3").Value = 101.287919427319
$ws.Range("E3").Value = 100.064826965332
$ws.Range("F3").Value = 107.347775421623
$ws.Range("G3").Value = 97.90897799573791
$ws.Range("H3").Value = 934735206
$ws.Range("I3").Value = "IBM"

$ws.Range("D4").Value = 90.50640963266871
$ws.Range("E4").Value = 87.2489013671875
$ws.Range("F4").Value = 95.38956045835909
$ws.Range("G4").Value = 85.53606518821447
$ws.Range("H4").Value = 934735206
$ws.Range("I4").Value = "IBM"

$ws.Range("D5").Value = 85.25074014806478
$ws.Range("E5").Value = 78.45457458496094
$ws.Range("F5").Value = 86.06175832639248
$ws.Range("G5").Value = 74.18575154447578
$ws.Range("H5").Value = 934735206
$ws.Range("I5").Value = "IBM"

$ws.Range("D6").Value = 95.59106736398724
$ws.Range("E6").Value = 92.6885986328125
$ws.Range("F6").Value = 97.50276889350413
$ws.Range("G6").Value = 90.57366306543928
$ws.Range("H6").Value = 934735206
$ws.Range("I6").Value = "IBM"

$ws.Range("D7").Value = 97.32796891671452
$ws.Range("E7").Value = 102.9965667724609
$ws.Range("F7").Value = 104.9074660205587
$ws.Range("G7").Value = 96.13525861626501
$ws.Range("H7").Value = 934735206
$ws.Range("I7").Value = "IBM"

$ws.Range("D8").Value = 102.2303641073759
$ws.Range("E8").Value = 99.40392303466795
$ws.Range("F8").Value = 102.5343481038552
$ws.Range("G8").Value = 95.58791066768551
$ws.Range("H8").Value = 934735206
$ws.Range("I8").Value = "IBM"

$ws.Range("D9").Value = 108.9924859214732
$ws.Range("E9").Value = 113.9004058837891
$ws.Range("F9").Value = 116.9874380293206
$ws.Range("G9").Value = 107.9090824623276
$ws.Range("H9").Value = 934735206
$ws.Range("I9").Value = "IBM"

$ws.Range("D10").Value = 114.3405407535148
$ws.Range("E10").Value = 105.4403762817383
$ws.Range("F10").Value = 115.9916509427199
$ws.Range("G10").Value = 104.9864894903359
$ws.Range("H10").Value = 934735206
$ws.Range("I10").Value = "IBM"

$ws.Range("D11").Value = 102.0133703046265
$ws.Range("E11").Value = 96.09503173828124
$ws.Range("F11").Value = 103.6407474942726
$ws.Range("G11").Value = 95.41086371730032
$ws.Range("H11").Value = 934735206
$ws.Range("I11").Value = "IBM"

$ws.Range("D12").Value = 97.5667578111136
$ws.Range("E12").Value = 103.4133834838867
$ws.Range("F12").Value = 109.0854833369854
$ws.Range("G12").Value = 97.47278813865412
$ws.Range("H12").Value = 934735206
$ws.Range("I12").Value = "IBM"

$ws.Range("D13").Value = 104.74532117344
$ws.Range("E13").Value = 110.9825820922852
$ws.Range("F13").Value = 116.019842037219
$ws.Range("G13").Value = 104.0944739212008
$ws.Range("H13").Value = 934735206
$ws.Range("I13").Value = "IBM"

$ws.Range("D14").Value = 104.9824737778915
$ws.Range("E14").Value = 99.24520111083984
$ws.Range("F14").Value = 110.9114441715889
$ws.Range("G14").Value = 98.66325651013037
$ws.Range("H14").Value = 934735206
$ws.Range("I14").Value = "IBM"

$ws.Range("D15").Value = 95.72274275629188
$ws.Range("E15").Value = 100.3261337280273
$ws.Range("F15").Value = 104.2095864325181
$ws.Range("G15").Value = 95.66736713904488
$ws.Range("H15").Value = 934735206
$ws.Range("I15").Value = "IBM"

$ws.Range("D16").Value = 106.1541138903501
$ws.Range("E16").Value = 80.76837158203125
$ws.Range("F16").Value = 108.0083660661066
$ws.Range("G16").Value = 79.83074759888396
$ws.Range("H16").Value = 934735206
$ws.Range("I16").Value = "IBM"

$ws.Range("D17").Value = 79.37354849812554
$ws.Range("E17").Value = 95.25392913818359
$ws.Range("F17").Value = 95.95546767475048
$ws.Range("G17").Value = 79.14679005229341
$ws.Range("H17").Value = 934735206
$ws.Range("I17").Value = "IBM"

$ws.Range("D18").Value = 101.4464652532417
$ws.Range("E18").Value = 100.5575256347656
$ws.Range("F18").Value = 104.2279786129962
$ws.Range("G18").Value = 97.6828172692513
$ws.Range("H18").Value = 934735206
$ws.Range("I18").Value = "IBM"

$ws.Range("D19").Value = 101.2659952316863
$ws.Range("E19").Value = 107.5334625244141
$ws.Range("F19").Value = 110.2174524525543
$ws.Range("G19").Value = 100.9250606642894
$ws.Range("H19").Value = 934735206
$ws.Range("I19").Value = "IBM"

$ws.Range("D20").Value = 106.8555202726466
$ws.Range("E20").Value = 98.15089416503906
$ws.Range("F20").Value = 108.1472703377638
$ws.Range("G20").Value = 96.07382177138196
$ws.Range("H20").Value = 934735206
$ws.Range("I20").Value = "IBM"

$ws.Range("D21").Value = 100.2532864445648
$ws.Range("E21").Value = 106.7363357543945
$ws.Range("F21").Value = 108.2661257886959
$ws.Range("G21").Value = 98.91657659078029
$ws.Range("H21").Value = 934735206
$ws.Range("I21").Value = "IBM"

$ws.Range("D22").Value = 79.80949509475371
$ws.Range("E22").Value = 94.21662139892578
$ws.Range("F22").Value = 97.0305174656789
$ws.Range("G22").Value = 78.42880841463099
$ws.Range("H22").Value = 934735206
$ws.Range("I22").Value = "IBM"

$ws.Range("D23").Value = 91.45749270661408
$ws.Range("E23").Value = 93.48785400390624
$ws.Range("F23").Value = 100.5066649121606
$ws.Range("G23").Value = 87.60208679654392
$ws.Range("H23").Value = 934735206
$ws.Range("I23").Value = "IBM"

$ws.Range("D24").Value = 94.265118011548
$ws.Range("E24").Value = 86.02192687988281
$ws.Range("F24").Value = 104.38806392411
$ws.Range("G24").Value = 81.5998789323984
$ws.Range("H24").Value = 934735206
$ws.Range("I24").Value = "IBM"

$ws.Range("D25").Value = 98.35965972325714
$ws.Range("E25").Value = 93.09192657470705
$ws.Range("F25").Value = 103.3538443528135
$ws.Range("G25").Value = 91.72419021687924
$ws.Range("H25").Value = 934735206
$ws.Range("I25").Value = "IBM"

$ws.Range("D26").Value = 105.9387992312513
$ws.Range("E26").Value = 112.3698883056641
$ws.Range("F26").Value = 117.8030546689078
$ws.Range("G26").Value = 103.2618138998987
$ws.Range("H26").Value = 934735206
$ws.Range("I26").Value = "IBM"

$ws.Range("D27").Value = 117.6937536962934
$ws.Range("E27").Value = 112.8886260986328
$ws.Range("F27").Value = 118.126214285328
$ws.Range("G27").Value = 109.0845502401273
$ws.Range("H27").Value = 934735206
$ws.Range("I27").Value = "IBM"

$ws.Range("D28").Value = 114.2207007437652
$ws.Range("E28").Value = 101.3404922485352
$ws.Range("F28").Value = 118.2710723167638
$ws.Range("G28").Value = 100.9516543019578
$ws.Range("H28").Value = 934735206
$ws.Range("I28").Value = "IBM"

$ws.Range("D29").Value = 115.1187305539174
$ws.Range("E29").Value = 114.6894073486328
$ws.Range("F29").Value = 122.0995169625494
$ws.Range("G29").Value = 106.635299824885
$ws.Range("H29").Value = 934735206
$ws.Range("I29").Value = "IBM"

$ws.Range("D30").Value = 112.6731659260076
$ws.Range("E30").Value = 114.8890914916992
$ws.Range("F30").Value = 123.2922171852099
$ws.Range("G30").Value = 108.5454663807936
$ws.Range("H30").Value = 934735206
$ws.Range("I30").Value = "IBM"

$ws.Range("D31").Value = 124.0138030856612
$ws.Range("E31").Value = 115.0337905883789
$ws.Range("F31").Value = 124.7789903420327
$ws.Range("G31").Value = 110.0556513439341
$ws.Range("H31").Value = 934735206
$ws.Range("I31").Value = "IBM"

$ws.Range("D32").Value = 107.015934426711
$ws.Range("E32").Value = 123.1627197265625
$ws.Range("F32").Value = 123.6703750575566
$ws.Range("G32").Value = 102.9102128222763
$ws.Range("H32").Value = 934735206
$ws.Range("I32").Value = "IBM"

$ws.Range("D33").Value = 127.1636705547864
$ws.Range("E33").Value = 121.4228210449219
$ws.Range("F33").Value = 132.6431416817111
$ws.Range("G33").Value = 119.845667155165
$ws.Range("H33").Value = 934735206
$ws.Range("I33").Value = "IBM"

$ws.Range("D34").Value = 119.4840333906759
$ws.Range("E34").Value = 115.3239440917969
$ws.Range("F34").Value = 120.9802061020377
$ws.Range("G34").Value = 113.6361820929791
$ws.Range("H34").Value = 934735206
$ws.Range("I34").Value = "IBM"

$ws.Range("D35").Value = 123.3788971423908
$ws.Range("E35").Value = 133.3291015625
$ws.Range("F35").Value = 133.7175050568041
$ws.Range("G35").Value = 121.6496366238282
$ws.Range("H35").Value = 934735206
$ws.Range("I35").Value = "IBM"

$ws.Range("D36").Value = 130.9909219225333
$ws.Range("E36").Value = 135.2936859130859
$ws.Range("F36").Value = 135.4059272140011
$ws.Range("G36").Value = 127.0903797149274
$ws.Range("H36").Value = 934735206
$ws.Range("I36").Value = "IBM"

$ws.Range("D37").Value = 154.0356052643627
$ws.Range("E37").Value = 173.7405853271484
$ws.Range("F37").Value = 186.2654879035313
$ws.Range("G37").Value = 149.3624107823105
$ws.Range("H37").Value = 934735206
$ws.Range("I37").Value = "IBM"

$ws.Range("D38").Value = 181.3768532612173
$ws.Range("E38").Value = 158.6570129394531
$ws.Range("F38").Value = 184.5079893522148
$ws.Range("G38").Value = 157.7596724927156
$ws.Range("H38").Value = 934735206
$ws.Range("I38").Value = "IBM"

$ws.Range("D39").Value = 167.2216258037211
$ws.Range("E39").Value = 185.2404937744141
$ws.Range("F39").Value = 189.212545050985
$ws.Range("G39").Value = 167.1541470652925
$ws.Range("H39").Value = 934735206
$ws.Range("I39").Value = "IBM"

$ws.Range("D40").Value = 214.5678431602001
$ws.Range("E40").Value = 201.0400390625
$ws.Range("F40").Value = 230.8478754297024
$ws.Range("G40").Value = 197.9182324093656
$ws.Range("H40").Value = 934735206
$ws.Range("I40").Value = "IBM"

$ws.Range("D41").Value = 217.425305098766
$ws.Range("E41").Value = 250.6340637207031
$ws.Range("F41").Value = 256.6132014302056
$ws.Range("G41").Value = 210.3581431757415
$ws.Range("H41").Value = 934735206
$ws.Range("I41").Value = "IBM"

$ws.Range("D42").Value = 244.7357116053205
$ws.Range("E42").Value = 238.6082000732422
$ws.Range("F42").Value = 249.4324847722812
$ws.Range("G42").Value = 211.6510518796289
$ws.Range("H42").Value = 934735206
$ws.Range("I42").Value = "IBM"

$ws.Range("D43").Value = 292.5718678669048
$ws.Range("E43").Value = 251.4499053955078
$ws.Range("F43").Value = 293.624746761845
$ws.Range("G43").Value = 250.5261583049219
$ws.Range("H43").Value = 934735206
$ws.Range("I43").Value = "IBM"

$ws.Range("D44").Value = 280.2000122070312
$ws.Range("E44").Value = 307.4599914550781
$ws.Range("F44").Value = 310.75
$ws.Range("G44").Value = 263.5599975585937
$ws.Range("H44").Value = 934735206
$ws.Range("I44").Value = "IBM"
